# Season 12, Matchday 11
# Add two new worksheets ("10" and "11") with scorer data for matchdays 10 and 11,
# mirroring the structure/style of the existing matchday sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet "10" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet10 = $wb.Worksheets.Add($null, $lastSheet)
$newSheet10.Name = "10"

$newSheet10.Range("B1").Value = "Scorer"
$newSheet10.Range("C1").Value = "Team"
$newSheet10.Range("D1").Value = "Goals"
$newSheet10.Range("E1").Value = "Picks"
$newSheet10.Range("F1").Value = "Matchday"

$newSheet10.Cells.Item(2, 2).Value = "агуеро"
$newSheet10.Cells.Item(3, 2).Value = "бензема"
$newSheet10.Cells.Item(4, 2).Value = "вернер"
$newSheet10.Cells.Item(5, 2).Value = "гнабри"
$newSheet10.Cells.Item(6, 2).Value = "госенс"
$newSheet10.Cells.Item(7, 2).Value = "гюндоган"
$newSheet10.Cells.Item(8, 2).Value = "де брюйне"
$newSheet10.Cells.Item(9, 2).Value = "жезус"
$newSheet10.Cells.Item(10, 2).Value = "жиру"
$newSheet10.Cells.Item(11, 2).Value = "иммобиле"
$newSheet10.Cells.Item(12, 2).Value = "карраско"
$newSheet10.Cells.Item(13, 2).Value = "кейн"
$newSheet10.Cells.Item(14, 2).Value = "коман"
$newSheet10.Cells.Item(15, 2).Value = "левандовски"
$newSheet10.Cells.Item(16, 2).Value = "ляказетт"
$newSheet10.Cells.Item(17, 2).Value = "марез"
$newSheet10.Cells.Item(18, 2).Value = "марсьяль"
$newSheet10.Cells.Item(19, 2).Value = "морелос"
$newSheet10.Cells.Item(20, 2).Value = "морено"
$newSheet10.Cells.Item(21, 2).Value = "муриэль"
$newSheet10.Cells.Item(22, 2).Value = "нерес"
$newSheet10.Cells.Item(23, 2).Value = "нсаме"
$newSheet10.Cells.Item(24, 2).Value = "обамеянг"
$newSheet10.Cells.Item(25, 2).Value = "пако"
$newSheet10.Cells.Item(26, 2).Value = "сапата"
$newSheet10.Cells.Item(27, 2).Value = "сон"
$newSheet10.Cells.Item(28, 2).Value = "стерлинг"
$newSheet10.Cells.Item(29, 2).Value = "суарес"
$newSheet10.Cells.Item(30, 2).Value = "тадич"
$newSheet10.Cells.Item(31, 2).Value = "фернандеш"
$newSheet10.Cells.Item(32, 2).Value = "эйкрейм"
$newSheet10.Cells.Item(33, 2).Value = "эрнандес"
$newSheet10.Cells.Item(2, 3).Value = "манчестер сити"
$newSheet10.Cells.Item(3, 3).Value = "реал мадрид"
$newSheet10.Cells.Item(4, 3).Value = "челси"
$newSheet10.Cells.Item(5, 3).Value = "бавария"
$newSheet10.Cells.Item(6, 3).Value = "аталанта"
$newSheet10.Cells.Item(7, 3).Value = "манчестер сити"
$newSheet10.Cells.Item(8, 3).Value = "манчестер сити"
$newSheet10.Cells.Item(9, 3).Value = "манчестер сити"
$newSheet10.Cells.Item(10, 3).Value = "челси"
$newSheet10.Cells.Item(11, 3).Value = "лацио"
$newSheet10.Cells.Item(12, 3).Value = "атлетико"
$newSheet10.Cells.Item(13, 3).Value = "тоттенхэм"
$newSheet10.Cells.Item(14, 3).Value = "бавария"
$newSheet10.Cells.Item(15, 3).Value = "бавария"
$newSheet10.Cells.Item(16, 3).Value = "арсенал"
$newSheet10.Cells.Item(17, 3).Value = "манчестер сити"
$newSheet10.Cells.Item(18, 3).Value = "манчестер юнайтед"
$newSheet10.Cells.Item(19, 3).Value = "рейнджерс"
$newSheet10.Cells.Item(20, 3).Value = "вильярреал"
$newSheet10.Cells.Item(21, 3).Value = "аталанта"
$newSheet10.Cells.Item(22, 3).Value = "аякс"
$newSheet10.Cells.Item(23, 3).Value = "янг бойз"
$newSheet10.Cells.Item(24, 3).Value = "арсенал"
$newSheet10.Cells.Item(25, 3).Value = "вильярреал"
$newSheet10.Cells.Item(26, 3).Value = "аталанта"
$newSheet10.Cells.Item(27, 3).Value = "тоттенхэм"
$newSheet10.Cells.Item(28, 3).Value = "сити"
$newSheet10.Cells.Item(29, 3).Value = "атлетико"
$newSheet10.Cells.Item(30, 3).Value = "аякс"
$newSheet10.Cells.Item(31, 3).Value = "манчестер юнайтед"
$newSheet10.Cells.Item(32, 3).Value = "мольде"
$newSheet10.Cells.Item(33, 3).Value = "милан"
$newSheet10.Cells.Item(2, 1).Value = 1
$newSheet10.Cells.Item(2, 4).Value = 0
$newSheet10.Cells.Item(2, 5).Value = 1
$newSheet10.Cells.Item(2, 6).Value = 10
$newSheet10.Cells.Item(3, 1).Value = 69
$newSheet10.Cells.Item(3, 4).Value = 1
$newSheet10.Cells.Item(3, 5).Value = 8
$newSheet10.Cells.Item(3, 6).Value = 10
$newSheet10.Cells.Item(4, 1).Value = 34
$newSheet10.Cells.Item(4, 4).Value = 0
$newSheet10.Cells.Item(4, 5).Value = 1
$newSheet10.Cells.Item(4, 6).Value = 10
$newSheet10.Cells.Item(5, 1).Value = 19
$newSheet10.Cells.Item(5, 4).Value = 0
$newSheet10.Cells.Item(5, 5).Value = 2
$newSheet10.Cells.Item(5, 6).Value = 10
$newSheet10.Cells.Item(6, 1).Value = 37
$newSheet10.Cells.Item(6, 4).Value = 0
$newSheet10.Cells.Item(6, 5).Value = 1
$newSheet10.Cells.Item(6, 6).Value = 10
$newSheet10.Cells.Item(7, 1).Value = 45
$newSheet10.Cells.Item(7, 4).Value = 1
$newSheet10.Cells.Item(7, 5).Value = 1
$newSheet10.Cells.Item(7, 6).Value = 10
$newSheet10.Cells.Item(8, 1).Value = 22
$newSheet10.Cells.Item(8, 4).Value = 1
$newSheet10.Cells.Item(8, 5).Value = 3
$newSheet10.Cells.Item(8, 6).Value = 10
$newSheet10.Cells.Item(9, 1).Value = 49
$newSheet10.Cells.Item(9, 4).Value = 0
$newSheet10.Cells.Item(9, 5).Value = 1
$newSheet10.Cells.Item(9, 6).Value = 10
$newSheet10.Cells.Item(10, 1).Value = 26
$newSheet10.Cells.Item(10, 4).Value = 0
$newSheet10.Cells.Item(10, 5).Value = 1
$newSheet10.Cells.Item(10, 6).Value = 10
$newSheet10.Cells.Item(11, 1).Value = 44
$newSheet10.Cells.Item(11, 4).Value = 0
$newSheet10.Cells.Item(11, 5).Value = 1
$newSheet10.Cells.Item(11, 6).Value = 10
$newSheet10.Cells.Item(12, 1).Value = 36
$newSheet10.Cells.Item(12, 4).Value = 0
$newSheet10.Cells.Item(12, 5).Value = 1
$newSheet10.Cells.Item(12, 6).Value = 10
$newSheet10.Cells.Item(13, 1).Value = 71
$newSheet10.Cells.Item(13, 4).Value = 0
$newSheet10.Cells.Item(13, 5).Value = 3
$newSheet10.Cells.Item(13, 6).Value = 10
$newSheet10.Cells.Item(14, 1).Value = 41
$newSheet10.Cells.Item(14, 4).Value = 0
$newSheet10.Cells.Item(14, 5).Value = 1
$newSheet10.Cells.Item(14, 6).Value = 10
$newSheet10.Cells.Item(15, 1).Value = 20
$newSheet10.Cells.Item(15, 4).Value = 1
$newSheet10.Cells.Item(15, 5).Value = 18
$newSheet10.Cells.Item(15, 6).Value = 10
$newSheet10.Cells.Item(16, 1).Value = 58
$newSheet10.Cells.Item(16, 4).Value = 0
$newSheet10.Cells.Item(16, 5).Value = 1
$newSheet10.Cells.Item(16, 6).Value = 10
$newSheet10.Cells.Item(17, 1).Value = 50
$newSheet10.Cells.Item(17, 4).Value = 0
$newSheet10.Cells.Item(17, 5).Value = 3
$newSheet10.Cells.Item(17, 6).Value = 10
$newSheet10.Cells.Item(18, 1).Value = 55
$newSheet10.Cells.Item(18, 4).Value = 0
$newSheet10.Cells.Item(18, 5).Value = 1
$newSheet10.Cells.Item(18, 6).Value = 10
$newSheet10.Cells.Item(19, 1).Value = 8
$newSheet10.Cells.Item(19, 4).Value = 0
$newSheet10.Cells.Item(19, 5).Value = 1
$newSheet10.Cells.Item(19, 6).Value = 10
$newSheet10.Cells.Item(20, 1).Value = 62
$newSheet10.Cells.Item(20, 4).Value = 2
$newSheet10.Cells.Item(20, 5).Value = 5
$newSheet10.Cells.Item(20, 6).Value = 10
$newSheet10.Cells.Item(21, 1).Value = 64
$newSheet10.Cells.Item(21, 4).Value = 1
$newSheet10.Cells.Item(21, 5).Value = 1
$newSheet10.Cells.Item(21, 6).Value = 10
$newSheet10.Cells.Item(22, 1).Value = 68
$newSheet10.Cells.Item(22, 4).Value = 1
$newSheet10.Cells.Item(22, 5).Value = 1
$newSheet10.Cells.Item(22, 6).Value = 10
$newSheet10.Cells.Item(23, 1).Value = 76
$newSheet10.Cells.Item(23, 4).Value = 0
$newSheet10.Cells.Item(23, 5).Value = 3
$newSheet10.Cells.Item(23, 6).Value = 10
$newSheet10.Cells.Item(24, 1).Value = 66
$newSheet10.Cells.Item(24, 4).Value = 0
$newSheet10.Cells.Item(24, 5).Value = 4
$newSheet10.Cells.Item(24, 6).Value = 10
$newSheet10.Cells.Item(25, 1).Value = 63
$newSheet10.Cells.Item(25, 4).Value = 0
$newSheet10.Cells.Item(25, 5).Value = 1
$newSheet10.Cells.Item(25, 6).Value = 10
$newSheet10.Cells.Item(26, 1).Value = 23
$newSheet10.Cells.Item(26, 4).Value = 0
$newSheet10.Cells.Item(26, 5).Value = 1
$newSheet10.Cells.Item(26, 6).Value = 10
$newSheet10.Cells.Item(27, 1).Value = 32
$newSheet10.Cells.Item(27, 4).Value = 0
$newSheet10.Cells.Item(27, 5).Value = 1
$newSheet10.Cells.Item(27, 6).Value = 10
$newSheet10.Cells.Item(28, 1).Value = 7
$newSheet10.Cells.Item(28, 4).Value = 0
$newSheet10.Cells.Item(28, 5).Value = 1
$newSheet10.Cells.Item(28, 6).Value = 10
$newSheet10.Cells.Item(29, 1).Value = 70
$newSheet10.Cells.Item(29, 4).Value = 0
$newSheet10.Cells.Item(29, 5).Value = 1
$newSheet10.Cells.Item(29, 6).Value = 10
$newSheet10.Cells.Item(30, 1).Value = 28
$newSheet10.Cells.Item(30, 4).Value = 1
$newSheet10.Cells.Item(30, 5).Value = 3
$newSheet10.Cells.Item(30, 6).Value = 10
$newSheet10.Cells.Item(31, 1).Value = 11
$newSheet10.Cells.Item(31, 4).Value = 0
$newSheet10.Cells.Item(31, 5).Value = 3
$newSheet10.Cells.Item(31, 6).Value = 10
$newSheet10.Cells.Item(32, 1).Value = 59
$newSheet10.Cells.Item(32, 4).Value = 0
$newSheet10.Cells.Item(32, 5).Value = 1
$newSheet10.Cells.Item(32, 6).Value = 10
$newSheet10.Cells.Item(33, 1).Value = 67
$newSheet10.Cells.Item(33, 4).Value = 0
$newSheet10.Cells.Item(33, 5).Value = 1
$newSheet10.Cells.Item(33, 6).Value = 10

$ws1.Range("B1:F1").Copy()
$newSheet10.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet10.Range("A2:A33").PasteSpecial(-4122)

# --- Sheet "11" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet11 = $wb.Worksheets.Add($null, $lastSheet)
$newSheet11.Name = "11"

$newSheet11.Range("B1").Value = "Scorer"
$newSheet11.Range("C1").Value = "Team"
$newSheet11.Range("D1").Value = "Goals"
$newSheet11.Range("E1").Value = "Picks"
$newSheet11.Range("F1").Value = "Matchday"

$newSheet11.Cells.Item(2, 2).Value = "абубакар"
$newSheet11.Cells.Item(3, 2).Value = "ауар"
$newSheet11.Cells.Item(4, 2).Value = "боаду"
$newSheet11.Cells.Item(5, 2).Value = "бэмфорд"
$newSheet11.Cells.Item(6, 2).Value = "депай"
$newSheet11.Cells.Item(7, 2).Value = "долберг"
$newSheet11.Cells.Item(8, 2).Value = "ибрагимович"
$newSheet11.Cells.Item(9, 2).Value = "инсинье"
$newSheet11.Cells.Item(10, 2).Value = "кессье"
$newSheet11.Cells.Item(11, 2).Value = "копмейнерс"
$newSheet11.Cells.Item(12, 2).Value = "корне"
$newSheet11.Cells.Item(13, 2).Value = "крессвелл"
$newSheet11.Cells.Item(14, 2).Value = "кюизанс"
$newSheet11.Cells.Item(15, 2).Value = "ларин"
$newSheet11.Cells.Item(16, 2).Value = "лингард"
$newSheet11.Cells.Item(17, 2).Value = "ляказетт"
$newSheet11.Cells.Item(18, 2).Value = "мален"
$newSheet11.Cells.Item(19, 2).Value = "мбаппе"
$newSheet11.Cells.Item(20, 2).Value = "мертенс"
$newSheet11.Cells.Item(21, 2).Value = "милик"
$newSheet11.Cells.Item(22, 2).Value = "морелос"
$newSheet11.Cells.Item(23, 2).Value = "неймар"
$newSheet11.Cells.Item(24, 2).Value = "осимхен"
$newSheet11.Cells.Item(25, 2).Value = "пайет"
$newSheet11.Cells.Item(26, 2).Value = "рафинья"
$newSheet11.Cells.Item(27, 2).Value = "сеферович"
$newSheet11.Cells.Item(28, 2).Value = "соучек"
$newSheet11.Cells.Item(29, 2).Value = "товен"
$newSheet11.Cells.Item(30, 2).Value = "эдегор"
$newSheet11.Cells.Item(31, 2).Value = "эдуард"
$newSheet11.Cells.Item(2, 3).Value = "бешикташ"
$newSheet11.Cells.Item(3, 3).Value = "лион"
$newSheet11.Cells.Item(4, 3).Value = "аз алкмар"
$newSheet11.Cells.Item(5, 3).Value = "лидс"
$newSheet11.Cells.Item(6, 3).Value = "лион"
$newSheet11.Cells.Item(7, 3).Value = "ницца"
$newSheet11.Cells.Item(8, 3).Value = "милан"
$newSheet11.Cells.Item(9, 3).Value = "наполи"
$newSheet11.Cells.Item(10, 3).Value = "милан"
$newSheet11.Cells.Item(11, 3).Value = "аз алкмар"
$newSheet11.Cells.Item(12, 3).Value = "лион"
$newSheet11.Cells.Item(13, 3).Value = "вэст хэм"
$newSheet11.Cells.Item(14, 3).Value = "марсель"
$newSheet11.Cells.Item(15, 3).Value = "бешикташ"
$newSheet11.Cells.Item(16, 3).Value = "вест хэм"
$newSheet11.Cells.Item(17, 3).Value = "арсенал"
$newSheet11.Cells.Item(18, 3).Value = "псв"
$newSheet11.Cells.Item(19, 3).Value = "псж"
$newSheet11.Cells.Item(20, 3).Value = "наполи"
$newSheet11.Cells.Item(21, 3).Value = "марсель"
$newSheet11.Cells.Item(22, 3).Value = "рейджерс"
$newSheet11.Cells.Item(23, 3).Value = "псж"
$newSheet11.Cells.Item(24, 3).Value = "наполи"
$newSheet11.Cells.Item(25, 3).Value = "марсель"
$newSheet11.Cells.Item(26, 3).Value = "лидс"
$newSheet11.Cells.Item(27, 3).Value = "бенфика"
$newSheet11.Cells.Item(28, 3).Value = "вест хэм"
$newSheet11.Cells.Item(29, 3).Value = "марсель"
$newSheet11.Cells.Item(30, 3).Value = "арсенал"
$newSheet11.Cells.Item(31, 3).Value = "селтик"
$newSheet11.Cells.Item(2, 1).Value = 65
$newSheet11.Cells.Item(2, 4).Value = 0
$newSheet11.Cells.Item(2, 5).Value = 8
$newSheet11.Cells.Item(2, 6).Value = 11
$newSheet11.Cells.Item(3, 1).Value = 67
$newSheet11.Cells.Item(3, 4).Value = 0
$newSheet11.Cells.Item(3, 5).Value = 1
$newSheet11.Cells.Item(3, 6).Value = 11
$newSheet11.Cells.Item(4, 1).Value = 35
$newSheet11.Cells.Item(4, 4).Value = 0
$newSheet11.Cells.Item(4, 5).Value = 4
$newSheet11.Cells.Item(4, 6).Value = 11
$newSheet11.Cells.Item(5, 1).Value = 34
$newSheet11.Cells.Item(5, 4).Value = 1
$newSheet11.Cells.Item(5, 5).Value = 1
$newSheet11.Cells.Item(5, 6).Value = 11
$newSheet11.Cells.Item(6, 1).Value = 10
$newSheet11.Cells.Item(6, 4).Value = 0
$newSheet11.Cells.Item(6, 5).Value = 2
$newSheet11.Cells.Item(6, 6).Value = 11
$newSheet11.Cells.Item(7, 1).Value = 1
$newSheet11.Cells.Item(7, 4).Value = 0
$newSheet11.Cells.Item(7, 5).Value = 1
$newSheet11.Cells.Item(7, 6).Value = 11
$newSheet11.Cells.Item(8, 1).Value = 28
$newSheet11.Cells.Item(8, 4).Value = 1
$newSheet11.Cells.Item(8, 5).Value = 8
$newSheet11.Cells.Item(8, 6).Value = 11
$newSheet11.Cells.Item(9, 1).Value = 29
$newSheet11.Cells.Item(9, 4).Value = 0
$newSheet11.Cells.Item(9, 5).Value = 2
$newSheet11.Cells.Item(9, 6).Value = 11
$newSheet11.Cells.Item(10, 1).Value = 40
$newSheet11.Cells.Item(10, 4).Value = 0
$newSheet11.Cells.Item(10, 5).Value = 1
$newSheet11.Cells.Item(10, 6).Value = 11
$newSheet11.Cells.Item(11, 1).Value = 57
$newSheet11.Cells.Item(11, 4).Value = 1
$newSheet11.Cells.Item(11, 5).Value = 3
$newSheet11.Cells.Item(11, 6).Value = 11
$newSheet11.Cells.Item(12, 1).Value = 76
$newSheet11.Cells.Item(12, 4).Value = 1
$newSheet11.Cells.Item(12, 5).Value = 1
$newSheet11.Cells.Item(12, 6).Value = 11
$newSheet11.Cells.Item(13, 1).Value = 55
$newSheet11.Cells.Item(13, 4).Value = 0
$newSheet11.Cells.Item(13, 5).Value = 1
$newSheet11.Cells.Item(13, 6).Value = 11
$newSheet11.Cells.Item(14, 1).Value = 68
$newSheet11.Cells.Item(14, 4).Value = 0
$newSheet11.Cells.Item(14, 5).Value = 1
$newSheet11.Cells.Item(14, 6).Value = 11
$newSheet11.Cells.Item(15, 1).Value = 37
$newSheet11.Cells.Item(15, 4).Value = 0
$newSheet11.Cells.Item(15, 5).Value = 3
$newSheet11.Cells.Item(15, 6).Value = 11
$newSheet11.Cells.Item(16, 1).Value = 24
$newSheet11.Cells.Item(16, 4).Value = 1
$newSheet11.Cells.Item(16, 5).Value = 2
$newSheet11.Cells.Item(16, 6).Value = 11
$newSheet11.Cells.Item(17, 1).Value = 18
$newSheet11.Cells.Item(17, 4).Value = 1
$newSheet11.Cells.Item(17, 5).Value = 1
$newSheet11.Cells.Item(17, 6).Value = 11
$newSheet11.Cells.Item(18, 1).Value = 61
$newSheet11.Cells.Item(18, 4).Value = 0
$newSheet11.Cells.Item(18, 5).Value = 3
$newSheet11.Cells.Item(18, 6).Value = 11
$newSheet11.Cells.Item(19, 1).Value = 53
$newSheet11.Cells.Item(19, 4).Value = 2
$newSheet11.Cells.Item(19, 5).Value = 12
$newSheet11.Cells.Item(19, 6).Value = 11
$newSheet11.Cells.Item(20, 1).Value = 32
$newSheet11.Cells.Item(20, 4).Value = 2
$newSheet11.Cells.Item(20, 5).Value = 1
$newSheet11.Cells.Item(20, 6).Value = 11
$newSheet11.Cells.Item(21, 1).Value = 7
$newSheet11.Cells.Item(21, 4).Value = 0
$newSheet11.Cells.Item(21, 5).Value = 1
$newSheet11.Cells.Item(21, 6).Value = 11
$newSheet11.Cells.Item(22, 1).Value = 20
$newSheet11.Cells.Item(22, 4).Value = 1
$newSheet11.Cells.Item(22, 5).Value = 1
$newSheet11.Cells.Item(22, 6).Value = 11
$newSheet11.Cells.Item(23, 1).Value = 19
$newSheet11.Cells.Item(23, 4).Value = 0
$newSheet11.Cells.Item(23, 5).Value = 1
$newSheet11.Cells.Item(23, 6).Value = 11
$newSheet11.Cells.Item(24, 1).Value = 60
$newSheet11.Cells.Item(24, 4).Value = 0
$newSheet11.Cells.Item(24, 5).Value = 1
$newSheet11.Cells.Item(24, 6).Value = 11
$newSheet11.Cells.Item(25, 1).Value = 54
$newSheet11.Cells.Item(25, 4).Value = 0
$newSheet11.Cells.Item(25, 5).Value = 1
$newSheet11.Cells.Item(25, 6).Value = 11
$newSheet11.Cells.Item(26, 1).Value = 63
$newSheet11.Cells.Item(26, 4).Value = 1
$newSheet11.Cells.Item(26, 5).Value = 3
$newSheet11.Cells.Item(26, 6).Value = 11
$newSheet11.Cells.Item(27, 1).Value = 52
$newSheet11.Cells.Item(27, 4).Value = 1
$newSheet11.Cells.Item(27, 5).Value = 4
$newSheet11.Cells.Item(27, 6).Value = 11
$newSheet11.Cells.Item(28, 1).Value = 30
$newSheet11.Cells.Item(28, 4).Value = 1
$newSheet11.Cells.Item(28, 5).Value = 1
$newSheet11.Cells.Item(28, 6).Value = 11
$newSheet11.Cells.Item(29, 1).Value = 0
$newSheet11.Cells.Item(29, 4).Value = 0
$newSheet11.Cells.Item(29, 5).Value = 2
$newSheet11.Cells.Item(29, 6).Value = 11
$newSheet11.Cells.Item(30, 1).Value = 42
$newSheet11.Cells.Item(30, 4).Value = 0
$newSheet11.Cells.Item(30, 5).Value = 2
$newSheet11.Cells.Item(30, 6).Value = 11
$newSheet11.Cells.Item(31, 1).Value = 41
$newSheet11.Cells.Item(31, 4).Value = 0
$newSheet11.Cells.Item(31, 5).Value = 2
$newSheet11.Cells.Item(31, 6).Value = 11

$ws1.Range("B1:F1").Copy()
$newSheet11.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet11.Range("A2:A31").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Restore the originally-active sheet/tab (adding sheets shifts activation to
# the newly created sheet by default).
$ws1.Activate()
